$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" field text (3/16/2012 -> 9/6/2012)
#    on the slide master and all three slide layouts.
# ---------------------------------------------------------------------------
function Set-DateText($container, $text) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $sh = $container.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $text
        }
    }
}

Set-DateText $p.SlideMaster "9/6/2012"
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    Set-DateText $p.SlideMaster.CustomLayouts.Item($i) "9/6/2012"
}

# ---------------------------------------------------------------------------
# 2) Slide 6 ("NOC Level" chart slide) geometry tweaks + two new shapes.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(6)

# Locate the two existing connectors we need to move/resize by name.
$connA = $null   # "Straight Connector 3" (diagonal connector, currently flipped)
$connB = $null   # "Straight Connector 4" (short horizontal connector)
for ($j = 1; $j -le $s.Shapes.Count; $j++) {
    $sh = $s.Shapes.Item($j)
    if ($sh.Name -eq "Straight Connector 3") { $connA = $sh }
    if ($sh.Name -eq "Straight Connector 4") { $connB = $sh }
}

# --- Straight Connector 3: un-flip, move down, shrink ---
$connA.VerticalFlip = 0
$connA.Left = 246
$connA.Top = 247.8709449
$connA.Width = 36
$connA.Height = 21.3790551

# --- Straight Connector 4: move down only (size unchanged) ---
$connB.Left = 282
$connB.Top = 269.2140945
$connB.Width = 45
$connB.Height = 0.0359055

# --- New connector: "Straight Connector 8" ---
$newConn = $s.Shapes.AddLine(493.35, 222.3750394, 493.35, 294.7568504)
$newConn.Width = 0
$newConn.VerticalFlip = -1
$newConn.Name = "Straight Connector 8"
$newConn.Line.BeginArrowheadStyle = 5   # msoArrowheadDiamond
$newConn.Line.BeginArrowheadWidth = 3   # msoArrowheadWide
$newConn.Line.BeginArrowheadLength = 3  # msoArrowheadLong

# --- New shape: "Rectangle 7" with note text ---
$newRect = $s.Shapes.AddShape(1, 370.05, 204, 124.05, 65.2140945)  # msoShapeRectangle
$newRect.Name = "Rectangle 7"
$newRect.Line.ForeColor.ObjectThemeColor = 1  # tx1/dk1 solid outline
$newRect.TextFrame.VerticalAnchor = 3         # msoAnchorMiddle -> anchor="ctr"

$tr = $newRect.TextFrame.TextRange
$tr.Text = "If comfortable death outcome is not achieved in first 3 days of hospitalization, expected LOS will double"
$tr.Font.Size = 10
$tr.Font.Bold = -1
$tr.Font.Color.ObjectThemeColor = 1
